$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Electric Rate Code (J2) and Electric Usage kWh (L2) were blank -> N/A
$ws.Range("J2").Value = "N/A"
$ws.Range("L2").Value = "N/A"

# Gas Usage (therms) P2: replace bad "6% 281.0" reading with the real usage
# value. Format the cell as Text first so the numeric-looking string isn't
# auto-coerced into a Number cell (matches the source data staying a string).
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "17342"
